# Update the metrics table:
#  - refresh RMSE/RRMSE/ME figures (columns B, D, F) for the existing
#    PAR / Temperature / Humidity / CO2 model rows (2-13)
#  - append three new "Leaf Temperature" model rows (14-16) for the
#    NN / GL / Combined models, matching the layout of the existing rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated metrics for existing rows -------------------------------

$ws.Range("B2").Value = 0.168116600438055
$ws.Range("D2").Value = 12.67880533766137
$ws.Range("F2").Value = -0.1668759703159332

$ws.Range("B3").Value = 2.207918313698289
$ws.Range("D3").Value = 166.5139934301305
$ws.Range("F3").Value = 1.491286757975844

$ws.Range("B4").Value = 0.1535421366873673
$ws.Range("D4").Value = 11.57964684698121
$ws.Range("F4").Value = -0.1504373073101044

$ws.Range("B5").Value = 1.171562587402149
$ws.Range("D5").Value = 5.285493133634397
$ws.Range("F5").Value = -1.13574721813202

$ws.Range("B6").Value = 2.017716160539228
$ws.Range("D6").Value = 9.102906687897267
$ws.Range("F6").Value = -1.659549614038649

$ws.Range("B7").Value = 1.1501706664674
$ws.Range("D7").Value = 5.18898369194372
$ws.Range("F7").Value = -1.098814272880555

$ws.Range("B8").Value = 6.626071919176557
$ws.Range("D8").Value = 12.97004535194824
$ws.Range("F8").Value = 6.546952104568479

$ws.Range("B9").Value = 46.37541056775313
$ws.Range("D9").Value = 90.77643370247736
$ws.Range("F9").Value = 46.25194408667794

$ws.Range("B10").Value = 6.626071919176557
$ws.Range("D10").Value = 12.97004535194824
$ws.Range("F10").Value = 6.546952104568479

$ws.Range("B11").Value = 1.951511303086819
$ws.Range("D11").Value = 0.4878778257717046
$ws.Range("F11").Value = -1.945152282714844

$ws.Range("B12").Value = 46.94930708714223
$ws.Range("D12").Value = 11.73732677178556
$ws.Range("F12").Value = 40.25923679430914

$ws.Range("B13").Value = 1.951511303086819
$ws.Range("D13").Value = 0.4878778257717046
$ws.Range("F13").Value = -1.945152282714844

# --- New "Leaf Temperature" rows --------------------------------------

$ws.Range("A14").Value = "Leaf Temperature (NN)"
$ws.Range("B14").Value = 1.237565324247811
$ws.Range("C14").Value = "°C"
$ws.Range("D14").Value = 5.772899471709905
$ws.Range("E14").Value = "%"
$ws.Range("F14").Value = -1.236618041992187
$ws.Range("G14").Value = "°C"

$ws.Range("A15").Value = "Leaf Temperature (GL)"
$ws.Range("B15").Value = 1.242150529272412
$ws.Range("C15").Value = "°C"
$ws.Range("D15").Value = 5.794288183194927
$ws.Range("E15").Value = "%"
$ws.Range("F15").Value = -0.7555795760566364
$ws.Range("G15").Value = "°C"

$ws.Range("A16").Value = "Leaf Temperature (Combined)"
$ws.Range("B16").Value = 0.9028122187080343
$ws.Range("C16").Value = "°C"
$ws.Range("D16").Value = 4.211368950241559
$ws.Range("E16").Value = "%"
$ws.Range("F16").Value = -0.6724491119384757
$ws.Range("G16").Value = "°C"
